# Rework the test fixture from a 3-column name/address table into a
# simple 2-column username/password table (supports the new "Read CSV
# file" / credential-reading helpers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old firstname/lastname/postcode grid (A1:C6) completely.
$ws.Range("A1:C6").ClearContents()

# Re-populate as a 2x2 username/password block. Written in this order so
# the shared-string table comes out as Password, Testusername,
# Testpassword, Username - matching how the sheet is actually consumed.
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "Testusername"
$ws.Range("B2").Value = "Testpassword"
$ws.Range("A1").Value = "Username"

# Widen the username column so the label/value are readable.
$ws.Columns.Item(1).ColumnWidth = 23.8

# Leave the selection on the natural top-left cell instead of the old E3.
$ws.Range("A1").Select()
